# Update gh-pages output data (regenerated "杭州-漫展信息" sheets):
# refreshed "want-to-go" counts, renamed one event, and swapped two cover
# image URLs, across the 展览 / 演出 / 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("C2").Value = "杭州·Aniidol偶像剧场09"
$ws.Range("F2").Value = 21
$ws.Range("F4").Value = 214
$ws.Range("F5").Value = 72
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 1304
$ws.Range("F10").Value = 1303
$ws.Range("F11").Value = 177
$ws.Range("F13").Value = 19
$ws.Range("F14").Value = 19
$ws.Range("F16").Value = 103
$ws.Range("F17").Value = 245
$ws.Range("F18").Value = 1659
$ws.Range("F19").Value = 613
$ws.Range("F22").Value = 2250
$ws.Range("F24").Value = 404
$ws.Range("F27").Value = 1211
$ws.Range("F30").Value = 2819
$ws.Range("F31").Value = 1617
$ws.Range("F34").Value = 663
$ws.Range("F35").Value = 865
$ws.Range("F36").Value = 1813
$ws.Range("F38").Value = 1829
$ws.Range("F42").Value = 39
$ws.Range("F43").Value = 861
$ws.Range("F44").Value = 792
$ws.Range("F45").Value = 1011
$ws.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202405/isG309e51715657222196.jpeg"
$ws.Range("F46").Value = 85
$ws.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202405/JbVl16OE1715676665714.jpeg"
$ws.Range("F49").Value = 3341

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 799
$ws.Range("F20").Value = 30

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("C2").Value = "杭州·Aniidol偶像剧场09"
$ws.Range("F2").Value = 21
$ws.Range("F4").Value = 214
$ws.Range("F5").Value = 72
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 1304
$ws.Range("F12").Value = 1303
$ws.Range("F13").Value = 177
$ws.Range("F15").Value = 19
$ws.Range("F16").Value = 19
$ws.Range("F18").Value = 103
$ws.Range("F19").Value = 245
$ws.Range("F20").Value = 1659
$ws.Range("F21").Value = 613
$ws.Range("F24").Value = 2250
$ws.Range("F25").Value = 404
$ws.Range("F27").Value = 1211
$ws.Range("F28").Value = 2819
$ws.Range("F29").Value = 1617
$ws.Range("F32").Value = 799
$ws.Range("F34").Value = 663
$ws.Range("F35").Value = 865
$ws.Range("F36").Value = 1813
$ws.Range("F39").Value = 1829
$ws.Range("F41").Value = 861
$ws.Range("F42").Value = 792
$ws.Range("F43").Value = 1011
$ws.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202405/isG309e51715657222196.jpeg"
$ws.Range("F44").Value = 85
$ws.Range("I44").Value = "//i2.hdslb.com/bfs/openplatform/202405/JbVl16OE1715676665714.jpeg"
$ws.Range("F46").Value = 30
$ws.Range("F48").Value = 3341
